$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 3.00002850027075
$ws.Range("B5").Value = 3.00002850027075
$ws.Range("B6").Value = 3.00002850027075
$ws.Range("B7").Value = 3.00002850027075
